$d = $word.ActiveDocument

# Locate the paragraph that currently holds "Соперникб кандидатб претендент"
# together with the _GoBack bookmark, and replace its contents with the
# post-edit structure: the original sentence (now carrying an English
# language tag on the paragraph mark), a new blank paragraph, and a new
# paragraph that keeps the bookmark and introduces the English phrase
# "In a similar vein".
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Соперникб кандидатб претендент*") {
        $targetPara = $candidate
    }
}

$xml = @'
<w:p w:rsidR="00B77EF8" w:rsidRPr="00B77EF8" w:rsidRDefault="00B77EF8" w:rsidP="00116A92" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Соперникб кандидатб претендент</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>In a similar vein</w:t></w:r></w:p>
'@

$targetPara.Range.InsertXML($xml)
